$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same serial date value (45205) for every
# data row (2 through 298). Bump it to 45206 for each of those rows, mirroring
# the upstream "Automatic update of files" commit.
for ($row = 2; $row -le 298; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
